$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.710.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.901.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5160"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.24%  "

$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07233"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8951"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07647"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.881.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.443"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008723"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9995"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "27.735.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.134"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.141.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.577"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.864"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.186"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.844"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08962"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.177"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("E33").Value = "  +0.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.817"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7794"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.615"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02084"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.054"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.090"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5498"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05270"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.662"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "113.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.500"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.41%  "

$ws.Range("E45").Value = "  -0.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4795"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9995"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.615"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.94%  "
